$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "43×39=" "87×32="
Replace-Text "84×28=" "26×36="
Replace-Text "27×81=" "16×79="
Replace-Text "32×98=" "34×94="
Replace-Text "56×22=" "54×60="
Replace-Text "73×27=" "95×79="
Replace-Text "15×89=" "80×73="
Replace-Text "78×78=" "12×64="
Replace-Text "19×68=" "71×59="
Replace-Text "15×69=" "97×42="
Replace-Text "65×28=" "37×99="
Replace-Text "92×38=" "33×69="
Replace-Text "14×53=" "26×24="
Replace-Text "60×76=" "59×47="
Replace-Text "48×21=" "39×74="
Replace-Text "46×44=" "76×61="
Replace-Text "49×84=" "80×79="
Replace-Text "56×96=" "54×69="
Replace-Text "43×35=" "91×48="
Replace-Text "22×95=" "19×25="
Replace-Text "38×89=" "81×34="
Replace-Text "44×70=" "41×14="
Replace-Text "64×49=" "77×56="
Replace-Text "43×64=" "54×50="
Replace-Text "14×17=" "58×50="

Write-Output "Done"
